$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date value in column C for rows 2-6 from 45233 to 45243
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45243
}
